$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.948.09"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.587.19"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "2.598.05"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.41"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.332"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("D14").Value = "3.042.91"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "58.908.26"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.46"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.598.42"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000133"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "345.91"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.42"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.21"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "0.0₃0717"
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.70"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.12"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.73"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.822"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("E40").Value = "  -4.85%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.77"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.596"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "267.77"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0515"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").Value = "1.955.22"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.10"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.35%  "
